# Update test data (placa numbers ABA070-073 -> ABA300-308, new extra rows of
# "procurador" results) and refresh the selection / active-sheet state, mirroring
# a re-run of the data-driven test with the "chrome driver actualizado".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("01-RegistrarAtencion")
$ws2 = $wb.Worksheets.Item("02-AsignarProcurador")

# ---------------------------------------------------------------------------
# Sheet 1: 01-RegistrarAtencion
# ---------------------------------------------------------------------------

# Rows 9 and 10 (and row 2) previously had no mark in column A; now they do.
$ws1.Range("A2").Value  = "X"
$ws1.Range("A9").Value  = "X"
$ws1.Range("A10").Value = "X"

# New plate numbers for each test row.
$ws1.Range("D2").Value  = "ABA300"
$ws1.Range("D3").Value  = "ABA301"
$ws1.Range("D4").Value  = "ABA302"
$ws1.Range("D5").Value  = "ABA303"
$ws1.Range("D6").Value  = "ABA304"
$ws1.Range("D7").Value  = "ABA305"
$ws1.Range("D8").Value  = "ABA306"
$ws1.Range("D9").Value  = "ABA307"
$ws1.Range("D10").Value = "ABA308"

# A couple of the "direccion" values shifted as the data cycled.
$ws1.Range("F6").Value  = "tambo real"
$ws1.Range("F7").Value  = "Miraflores"
$ws1.Range("F10").Value = "Miraflores"

# Column I (index 9) now has its own best-fit width instead of sharing the
# generic 9..1017 run.
$ws1.Columns.Item(9).ColumnWidth = 8.61

# ---------------------------------------------------------------------------
# Sheet 2: 02-AsignarProcurador
# ---------------------------------------------------------------------------

# New row 6, matching the layout/formula pattern of rows 2-5 (added first,
# ahead of reassigning the "procurador" names below, to mirror the order the
# shared strings were appended in the original session).
$ws2.Range("A5:F5").Copy()
$ws2.Range("A6").PasteSpecial(-4122)
$ws2.Range("A6").Value = "x"
$ws2.Range("B6").Value = "clegall@gmail.com"
$ws2.Range("C6").Value = "Test2019#"
$ws2.Range("D6").Formula = "='01-RegistrarAtencion'!D6"
$ws2.Range("E6").Value = "PROCURADOR AUTO5"
$ws2.Range("F6").Value = "En camino"

# Re-assigned "procurador" for the existing rows.
$ws2.Range("E2").Value = "PROCURADOR AUTO1"
$ws2.Range("E3").Value = "PROCURADOR AUTO2"
$ws2.Range("E4").Value = "PROCURADOR AUTO3"
$ws2.Range("E5").Value = "PROCURADOR AUTO4"

$ws2.Hyperlinks.Add($ws2.Range("B6"), "mailto:clegall@gmail.com")

# Adding the hyperlink re-stamps B6 with a generic "followed hyperlink"
# style; re-copy the sibling cell's format so it matches rows 2-5 exactly.
$ws2.Range("B5").Copy()
$ws2.Range("B6").PasteSpecial(-4122)

# Extra column widths now auto-sized for the new data.
$ws2.Columns.Item(2).ColumnWidth = 16.28
$ws2.Columns.Item(3).ColumnWidth = 9.28
$ws2.Columns.Item(4).ColumnWidth = 7.17

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping (matches what Excel records after
# re-running the data-driven test and leaving the cursor on sheet 2).
# ---------------------------------------------------------------------------

$ws1.Activate()
$ws1.Range("F11").Select()

$ws2.Activate()
$ws2.Range("I5").Select()
